$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the CELEBREX 200MG 15 CAPS. row (row 17) entirely - this shifts all
# subsequent rows up by one, Excel auto-recalculates the SUM total, and
# unused shared strings / merged cells are cleaned up automatically.
$ws.Rows.Item(17).Delete()

# Update the generated-timestamp cell (now at A57 after the row shift) to
# reflect the new save time.
$ws.Range("A57").Value = "Tuesday, 10 June, 2025 4:34 PM"
